$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the student ID in row 2 (231994 -> 231995), keeping it text
# (these IDs are stored as text, per the sheet's numberStoredAsText flag)
$ws.Range("A2").Value = "'231995"

# The former row 3 held a duplicate record (same ID, now merged into row 2);
# remove it so the log only has the single, updated row.
$ws.Rows("3").Delete()
